$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "248.99"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.99"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.356"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05618"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.408"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9263"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1452"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07496"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.03212"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.03090"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.09313"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.555"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.001605"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.04730"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0005761"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.006372"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.005065"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.001033"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.730"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.155"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.3307"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1319"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0003000"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.03943"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.006821"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1068"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.007751"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00005573"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6801"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1944"
$cell.Style = "Normal"
